$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status text for all data rows (E2:E11) from OPTIMAL to TIME_LIMIT
$ws.Range("E2:E11").Value = "TIME_LIMIT"

# Corrected fixed-recourse data: objective (B), gap (C), solve time (D)
$ws.Range("B2").Value = -493.878975631448
$ws.Range("C2").Value = 2.9786947810041657
$ws.Range("D2").Value = 3898.835644505

$ws.Range("B3").Value = -493.7770132109689
$ws.Range("C3").Value = 1.6976536168606413
$ws.Range("D3").Value = 4007.940096425

$ws.Range("B4").Value = -494.6697799631145
$ws.Range("C4").Value = 1.1096250214684438
$ws.Range("D4").Value = 3883.604352099

$ws.Range("B5").Value = -501.03532061828355
$ws.Range("C5").Value = 5.723580022259411
$ws.Range("D5").Value = 3634.586522726

$ws.Range("B6").Value = -489.92059559647726
$ws.Range("C6").Value = 5.000728538711415
$ws.Range("D6").Value = 3626.835960805

$ws.Range("B7").Value = -485.41422872748444
$ws.Range("C7").Value = 1.3593148623127218
$ws.Range("D7").Value = 3771.509533898

$ws.Range("B8").Value = -482.08357617734873
$ws.Range("C8").Value = 1.0366550606859448
$ws.Range("D8").Value = 3862.965649845

$ws.Range("B9").Value = -494.2150625290095
$ws.Range("C9").Value = 7.134439581090576
$ws.Range("D9").Value = 3603.486274451

$ws.Range("B10").Value = -490.8321604061183
$ws.Range("C10").Value = 1.3523059005132159
$ws.Range("D10").Value = 3720.388811662

$ws.Range("B11").Value = -483.7629393419096
$ws.Range("C11").Value = 1.1420716179603556
$ws.Range("D11").Value = 3846.992365969
